# Add "Expected Result" column to the Reqs importer test sheet.
#
# Before:  A=Sequence  B=Text  C=Linked Requirements  D=Notes
# After:   A=Sequence  B=Text  C=Linked Requirements  D=Expected Result  E=Notes
#
# Rows that had a "Linked Requirements"/"Notes" value keep "Notes" (now
# shifted to column E) and get no "Expected Result" value.
# Rows that had neither get "Yes" in the new "Expected Result" column (D).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = 21

# Header row: rename the old "Notes" header to "Expected Result" and add a
# new "Notes" header in the new column E.
$ws.Cells.Item(1, 4).Value = "Expected Result"
$ws.Cells.Item(1, 5).Value = "Notes"

for ($r = 2; $r -le $lastRow; $r++) {
    $linked = $ws.Cells.Item($r, 3).Value2
    if ($linked -eq "SRS-SW-0001") {
        # Row already has Linked Requirements / Notes -> shift Notes to E,
        # leave the new Expected Result column (D) blank.
        $ws.Cells.Item($r, 5).Value = "Notes"
        $ws.Cells.Item($r, 4).Value = $null
    } else {
        # Row has no Linked Requirements / Notes -> set Expected Result.
        $ws.Cells.Item($r, 4).Value = "Yes"
    }
}

# Widen the new Expected Result column (closest value Excel's pixel-snapped
# column-width model can represent to the target 20.140625 characters).
$ws.Columns.Item(4).ColumnWidth = 19.25

# Update the used-range dimension/selection bookkeeping to match the edit.
$ws.Range("G15").Select()
